# Update odds values in sheet1 (Jogos da Semana FlashScore 2025-04-17)
# per commit diff: values for several match rows were refreshed with updated odds.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.38
$ws.Range("I2").Value = 3.25
$ws.Range("J2").Value = 1.05
$ws.Range("L2").Value = 1.33
$ws.Range("M2").Value = 3.2
$ws.Range("N2").Value = 2.15
$ws.Range("O2").Value = 1.67
$ws.Range("R2").Value = 1.91
$ws.Range("S2").Value = 1.91
$ws.Range("T2").Value = 7.5
$ws.Range("U2").Value = 11
$ws.Range("W2").Value = 23
$ws.Range("X2").Value = 21
$ws.Range("Z2").Value = 8
$ws.Range("AA2").Value = 6
$ws.Range("AB2").Value = 15
$ws.Range("AD2").Value = 301
$ws.Range("AE2").Value = 9

# Row 3
$ws.Range("H3").Value = 3.25
$ws.Range("I3").Value = 5
$ws.Range("N3").Value = 2.6
$ws.Range("O3").Value = 1.48
$ws.Range("P3").Value = 1.57
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 2.38
$ws.Range("S3").Value = 1.53
$ws.Range("T3").Value = 5
$ws.Range("Z3").Value = 6.5
$ws.Range("AB3").Value = 23
$ws.Range("AC3").Value = 101
$ws.Range("AE3").Value = 9.5
$ws.Range("AG3").Value = 19
$ws.Range("AI3").Value = 51

# Row 4
$ws.Range("G4").Value = 2.38
$ws.Range("I4").Value = 3.3
$ws.Range("J4").Value = 1.07
$ws.Range("L4").Value = 1.41
$ws.Range("M4").Value = 2.62
$ws.Range("N4").Value = 2.5
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 1.57
$ws.Range("Q4").Value = 2.25
$ws.Range("R4").Value = 2.1
$ws.Range("S4").Value = 1.67
$ws.Range("W4").Value = 23
$ws.Range("X4").Value = 23
$ws.Range("Y4").Value = 41
$ws.Range("Z4").Value = 6.5
$ws.Range("AE4").Value = 8
$ws.Range("AI4").Value = 29

# Row 5
$ws.Range("J5").Value = 1.05
$ws.Range("K5").Value = 9
$ws.Range("L5").Value = 1.33

# Row 6
$ws.Range("G6").Value = 1.67
$ws.Range("H6").Value = 3.75
$ws.Range("J6").Value = 1.02
$ws.Range("L6").Value = 1.19
$ws.Range("M6").Value = 3.75
$ws.Range("N6").Value = 1.85
$ws.Range("O6").Value = 2
$ws.Range("P6").Value = 1.36
$ws.Range("Q6").Value = 3
$ws.Range("R6").Value = 1.8
$ws.Range("S6").Value = 1.91
$ws.Range("T6").Value = 7.5
$ws.Range("Z6").Value = 11
$ws.Range("AE6").Value = 15
$ws.Range("AG6").Value = 17

# Row 7
$ws.Range("G7").Value = 2.55
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 2.9
$ws.Range("J7").Value = 1.05
$ws.Range("K7").Value = 6.5
$ws.Range("M7").Value = 2.5
$ws.Range("W7").Value = 26
$ws.Range("Z7").Value = 6.5
$ws.Range("AE7").Value = 7

# Row 9
$ws.Range("G9").Value = 1.33
$ws.Range("I9").Value = 9
$ws.Range("N9").Value = 1.75
$ws.Range("O9").Value = 2.05
$ws.Range("U9").Value = 6.5
$ws.Range("W9").Value = 8.5
$ws.Range("AD9").Value = 1000
$ws.Range("AE9").Value = 19
$ws.Range("AG9").Value = 23
$ws.Range("AH9").Value = 101
$ws.Range("AI9").Value = 51
$ws.Range("AJ9").Value = 51

# Row 11
$ws.Range("J11").Value = 1.07
$ws.Range("K11").Value = 7
$ws.Range("L11").Value = 1.41
$ws.Range("M11").Value = 2.62

# Row 12
$ws.Range("G12").Value = 1.55
$ws.Range("H12").Value = 3.65
$ws.Range("I12").Value = 5.9
$ws.Range("L12").Value = 1.4
$ws.Range("M12").Value = 2.52
$ws.Range("R12").Value = 2.25
$ws.Range("U12").Value = 6
$ws.Range("V12").Value = 8.75
$ws.Range("W12").Value = 10.25
$ws.Range("Y12").Value = 40
$ws.Range("Z12").Value = 7.6
$ws.Range("AA12").Value = 7.4
$ws.Range("AB12").Value = 25
$ws.Range("AE12").Value = 11.75
$ws.Range("AF12").Value = 32
$ws.Range("AG12").Value = 21
$ws.Range("AI12").Value = 80
$ws.Range("AJ12").Value = 100

# Row 14
$ws.Range("G14").Value = 2.27
$ws.Range("H14").Value = 3.1
$ws.Range("I14").Value = 3.1
$ws.Range("L14").Value = 1.3
$ws.Range("M14").Value = 2.92
$ws.Range("N14").Value = 1.88
$ws.Range("O14").Value = 1.72
$ws.Range("P14").Value = 1.39
$ws.Range("Q14").Value = 2.55
$ws.Range("S14").Value = 1.93
$ws.Range("U14").Value = 11.25
$ws.Range("V14").Value = 8.75
$ws.Range("W14").Value = 23
$ws.Range("Y14").Value = 28
$ws.Range("Z14").Value = 9
$ws.Range("AA14").Value = 6
$ws.Range("AB14").Value = 13
$ws.Range("AC14").Value = 60
$ws.Range("AD14").Value = 450
$ws.Range("AE14").Value = 9.5
$ws.Range("AF14").Value = 16.5
$ws.Range("AG14").Value = 10.75
$ws.Range("AH14").Value = 40
$ws.Range("AI14").Value = 27

# Row 15
$ws.Range("G15").Value = 4.15
$ws.Range("H15").Value = 3.85
$ws.Range("I15").Value = 1.7
$ws.Range("N15").Value = 1.47
$ws.Range("O15").Value = 2.35
$ws.Range("R15").Value = 1.47
$ws.Range("S15").Value = 2.32
$ws.Range("T15").Value = 17.5
$ws.Range("U15").Value = 29
$ws.Range("V15").Value = 14
$ws.Range("W15").Value = 70
$ws.Range("Z15").Value = 16.5
$ws.Range("AA15").Value = 8.25
$ws.Range("AB15").Value = 12.5
$ws.Range("AD15").Value = 200
$ws.Range("AE15").Value = 10.5
$ws.Range("AF15").Value = 10.25
$ws.Range("AH15").Value = 15
$ws.Range("AI15").Value = 11.75
$ws.Range("AJ15").Value = 18

# Row 19
$ws.Range("G19").Value = 1.95
$ws.Range("H19").Value = 3.5
$ws.Range("I19").Value = 3.5
$ws.Range("N19").Value = 1.8
$ws.Range("O19").Value = 2
$ws.Range("P19").Value = 1.33
$ws.Range("Q19").Value = 3.25
$ws.Range("U19").Value = 10
$ws.Range("W19").Value = 17
$ws.Range("X19").Value = 15
$ws.Range("AE19").Value = 13
$ws.Range("AG19").Value = 13
$ws.Range("AH19").Value = 41

# Row 20
$ws.Range("N20").Value = 1.22
$ws.Range("O20").Value = 4

# Row 21
$ws.Range("G21").Value = 1.62
$ws.Range("I21").Value = 4.5
$ws.Range("J21").Value = 17
$ws.Range("K21").Value = 1.03
$ws.Range("L21").Value = 1.11
$ws.Range("M21").Value = 5
$ws.Range("N21").Value = 1.5
$ws.Range("O21").Value = 2.5
$ws.Range("P21").Value = 1.25
$ws.Range("Q21").Value = 3.75
$ws.Range("R21").Value = 1.57
$ws.Range("S21").Value = 2.25
$ws.Range("T21").Value = 10
$ws.Range("U21").Value = 10
$ws.Range("V21").Value = 9
$ws.Range("W21").Value = 13
$ws.Range("Y21").Value = 19
$ws.Range("Z21").Value = 17
$ws.Range("AA21").Value = 8.5
$ws.Range("AB21").Value = 13
$ws.Range("AD21").Value = 126
$ws.Range("AE21").Value = 19
$ws.Range("AG21").Value = 15
$ws.Range("AJ21").Value = 34

# Row 22
$ws.Range("J22").Value = 1.01
$ws.Range("L22").Value = 1.08

# Row 23
$ws.Range("G23").Value = 2.4
$ws.Range("I23").Value = 3
$ws.Range("J23").Value = 1.11
$ws.Range("K23").Value = 6
$ws.Range("L23").Value = 1.5
$ws.Range("M23").Value = 2.25
$ws.Range("N23").Value = 2.45
$ws.Range("O23").Value = 1.42
$ws.Range("P23").Value = 1.53
$ws.Range("Q23").Value = 2.18
$ws.Range("R23").Value = 2.05
$ws.Range("S23").Value = 1.6
$ws.Range("T23").Value = 6
$ws.Range("U23").Value = 10.25
$ws.Range("V23").Value = 10
$ws.Range("W23").Value = 25
$ws.Range("X23").Value = 25
$ws.Range("Y23").Value = 45
$ws.Range("Z23").Value = 6.3
$ws.Range("AB23").Value = 18.5
$ws.Range("AC23").Value = 120
$ws.Range("AE23").Value = 6.9
$ws.Range("AF23").Value = 13.5
$ws.Range("AH23").Value = 40
$ws.Range("AJ23").Value = 55
